$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 11 <-> 12: swap Id / Ost / Nord values ---
$ws.Range("A11").Value = 131066883
$ws.Range("Q11").Value = 408727
$ws.Range("R11").Value = 6703044

$ws.Range("A12").Value = 131066880
$ws.Range("Q12").Value = 408732
$ws.Range("R12").Value = 6703060

# --- Rows 13 <-> 14: swap full record (species differs) ---
$ws.Range("A13").Value = 131066896
$ws.Range("B13").Value = 79000
$ws.Range("E13").Value = 6446
$ws.Range("F13").Value = "Kolflarnlav"
$ws.Range("G13").Value = "Carbonicola anthracophila"
$ws.Range("H13").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q13").Value = 408691
$ws.Range("R13").Value = 6703020

$ws.Range("A14").Value = 131066888
$ws.Range("B14").Value = 91758
$ws.Range("E14").Value = 112
$ws.Range("F14").Value = "Stjärntagging"
$ws.Range("G14").Value = "Asterodon ferruginosus"
$ws.Range("H14").Value = "Pat."
$ws.Range("Q14").Value = 408773
$ws.Range("R14").Value = 6703124

# --- Rows 18 <-> 19: full row swap, including optional cells K,L,M,N,AC ---
$ws.Range("A18").Value = 131063926
$ws.Range("B18").Value = 83089
$ws.Range("E18").Value = 1312
$ws.Range("F18").Value = "Gammelgransskål"
$ws.Range("G18").Value = "Pseudographis pinicola"
$ws.Range("H18").Value = "(Nyl.) Rehm"
$ws.Range("K18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("P18").Value = "Torsby kommun, Vrm"
$ws.Range("Q18").Value = 408603
$ws.Range("R18").Value = 6702927
$ws.Range("S18").Value = 5
$ws.Range("AC18").ClearContents()
$ws.Range("AW18").Value = "Joakim Karlsson"
$ws.Range("AX18").Value = "Joakim Karlsson"

$ws.Range("A19").Value = 131066881
$ws.Range("B19").Value = 57884
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = "Tretåig hackspett"
$ws.Range("G19").Value = "Picoides tridactylus"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("K19").Value = ""
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = "äldre spår"
$ws.Range("N19").Value = ""
$ws.Range("P19").Value = "Färntjärnen, Vrm"
$ws.Range("Q19").Value = 408720
$ws.Range("R19").Value = 6703065
$ws.Range("S19").Value = 20
$ws.Range("AC19").Value = "Ringhack på gran"
$ws.Range("AW19").Value = "Moa Björnberg dillner"
$ws.Range("AX19").Value = "Moa Björnberg dillner"

# --- Rows 23 <-> 24: swap full record (species differs) ---
$ws.Range("A23").Value = 131066899
$ws.Range("B23").Value = 83089
$ws.Range("E23").Value = 1312
$ws.Range("F23").Value = "Gammelgransskål"
$ws.Range("G23").Value = "Pseudographis pinicola"
$ws.Range("H23").Value = "(Nyl.) Rehm"
$ws.Range("Q23").Value = 408724
$ws.Range("R23").Value = 6703049

$ws.Range("A24").Value = 131066877
$ws.Range("B24").Value = 57884
$ws.Range("E24").Value = 100109
$ws.Range("F24").Value = "Tretåig hackspett"
$ws.Range("G24").Value = "Picoides tridactylus"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("Q24").Value = 408778
$ws.Range("R24").Value = 6703117
